# Add season-record columns (Wins / Losses / Ties) to the roster table.
# Mirrors the existing "Unnamed: 28" header at AC1 so the new headers
# pick up the same bold / centered / bordered style, then fills the
# won-loss-tie record (71-91-0) down every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): AD1=Wins, AE1=Losses, AF1=Ties -------------------
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the formatting of the last existing header cell (AC1) onto the new
# header cells so they match the rest of the header row (bold, centered,
# bordered) instead of getting default formatting.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows (2-43): season record repeated for every player -----------
$lastRow = 43
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 71
    $ws.Cells.Item($r, 31).Value = 91
    $ws.Cells.Item($r, 32).Value = 0
}
